$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 472.3889
$ws.Range("I33").Value = 392.85715
$ws.Range("J33").Value = 750.75
$ws.Range("K33").Value = 392.85715
$ws.Range("L33").Value = 750.75
$ws.Range("M33").Value = -163.85715
$ws.Range("N33").Value = -1208.75

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2313.4285
$ws.Range("I62").Value = 2313.4285
$ws.Range("K62").Value = 2313.4285
$ws.Range("M62").Value = -1689.4285

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3512.9565
$ws.Range("I64").Value = 3328.5
$ws.Range("J64").Value = 3799.889
$ws.Range("K64").Value = 3328.5
$ws.Range("L64").Value = 3799.889
$ws.Range("M64").Value = -3080.5
$ws.Range("N64").Value = -4295.889

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2313.4285
$ws.Range("I65").Value = 2313.4285
$ws.Range("K65").Value = 11567.1425
$ws.Range("M65").Value = -8447.1425

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3512.9565
$ws.Range("I67").Value = 3328.5
$ws.Range("J67").Value = 3799.889
$ws.Range("K67").Value = 3328.5
$ws.Range("L67").Value = 3799.889
$ws.Range("M67").Value = -2470.5
$ws.Range("N67").Value = -5515.889

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1879.3846
$ws.Range("I100").Value = 1338.2667
$ws.Range("J100").Value = 2617.2727
$ws.Range("K100").Value = 1338.2667
$ws.Range("L100").Value = 2617.2727
$ws.Range("M100").Value = -797.2666999999999
$ws.Range("N100").Value = -3699.2727

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1766.1666
$ws.Range("I127").Value = 500
$ws.Range("J127").Value = 2399.25
$ws.Range("K127").Value = 1500
$ws.Range("L127").Value = 7197.75
$ws.Range("M127").Value = 3460
$ws.Range("N127").Value = -17117.75

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 969.0213
$ws.Range("J129").Value = 998.7954999999999
$ws.Range("L129").Value = 2996.3865
$ws.Range("N129").Value = -12996.3865

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 56991.05
$ws.Range("I132").Value = 67176.875
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 201530.625
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -199000.625
$ws.Range("N132").Value = -13060.0001

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5627.579
$ws.Range("I61").Value = 6845.778
$ws.Range("J61").Value = 4531.2
$ws.Range("K61").Value = 6845.778
$ws.Range("L61").Value = 4531.2
$ws.Range("M61").Value = -6633.778
$ws.Range("N61").Value = -4955.2

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 754.3214
$ws.Range("I74").Value = 337.70587
$ws.Range("J74").Value = 1398.1818
$ws.Range("K74").Value = 337.70587
$ws.Range("L74").Value = 1398.1818
$ws.Range("M74").Value = 536.29413
$ws.Range("N74").Value = -3146.1818

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 754.3214
$ws.Range("I77").Value = 337.70587
$ws.Range("J77").Value = 1398.1818
$ws.Range("K77").Value = 1688.52935
$ws.Range("L77").Value = 6990.909000000001
$ws.Range("M77").Value = 2679.47065
$ws.Range("N77").Value = -15726.909

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1838.4783
$ws.Range("I122").Value = 1811
$ws.Range("J122").Value = 1937.4
$ws.Range("K122").Value = 5433
$ws.Range("L122").Value = 5812.200000000001
$ws.Range("M122").Value = -2983
$ws.Range("N122").Value = -10712.2

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 13378.429
$ws.Range("I132").Value = 980.4516
$ws.Range("K132").Value = 2941.3548
$ws.Range("M132").Value = -411.3548000000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5627.579
$ws.Range("I136").Value = 6845.778
$ws.Range("J136").Value = 4531.2
$ws.Range("K136").Value = 20537.334
$ws.Range("L136").Value = 13593.6
$ws.Range("M136").Value = -17987.334
$ws.Range("N136").Value = -18693.6

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 991.55554
$ws.Range("I20").Value = 967.6
$ws.Range("J20").Value = 1021.5
$ws.Range("K20").Value = 967.6
$ws.Range("L20").Value = 1021.5
$ws.Range("M20").Value = -720.6
$ws.Range("N20").Value = -1515.5

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4688.778
$ws.Range("I62").Value = 4180
$ws.Range("K62").Value = 4180
$ws.Range("M62").Value = -3556

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4688.778
$ws.Range("I65").Value = 4180
$ws.Range("K65").Value = 20900
$ws.Range("M65").Value = -17780

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2856.348
$ws.Range("I132").Value = 805.7143
$ws.Range("K132").Value = 2417.1429
$ws.Range("M132").Value = 112.8571000000002

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 40
$ws.Range("K4").Value = 120
$ws.Range("M4").Value = -8

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 500000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1105.5
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 1167.8846
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3503.6538
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -5125.6538

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1105.5
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 1167.8846
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 10510.9614
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -18622.9614

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3080.6667
$ws.Range("J107").Value = 670.93335
$ws.Range("L107").Value = 2012.80005
$ws.Range("N107").Value = -5852.80005

# CUL row 119
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 3215
$ws.Range("I119").Value = 953.3333
$ws.Range("J119").Value = 10000
$ws.Range("K119").Value = 2859.9999
$ws.Range("L119").Value = 30000
$ws.Range("M119").Value = 1978.0001
$ws.Range("N119").Value = -39676

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 805.0700000000001
$ws.Range("J131").Value = 827.02106
$ws.Range("L131").Value = 2481.06318
$ws.Range("N131").Value = -12561.06318

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7816536
$ws.Range("I70").Value = 4333.3335
$ws.Range("K70").Value = 4333.3335
$ws.Range("M70").Value = -4063.3335

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7816536
$ws.Range("I73").Value = 4333.3335
$ws.Range("K73").Value = 4333.3335
$ws.Range("M73").Value = -3397.3335

# GSM row 124
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 48980
$ws.Range("J124").Value = 48980
$ws.Range("L124").Value = 48980
$ws.Range("N124").Value = -58800

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5600
$ws.Range("I126").Value = 4263.636
$ws.Range("J126").Value = 7437.5
$ws.Range("K126").Value = 12790.908
$ws.Range("L126").Value = 22312.5
$ws.Range("M126").Value = -10320.908
$ws.Range("N126").Value = -27252.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3511.1538
$ws.Range("I7").Value = 3714.1052
$ws.Range("J7").Value = 2960.2856
$ws.Range("K7").Value = 3714.1052
$ws.Range("L7").Value = 2960.2856
$ws.Range("M7").Value = -3602.1052
$ws.Range("N7").Value = -3184.2856

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 787
$ws.Range("I46").Value = 654.2308
$ws.Range("J46").Value = 1650
$ws.Range("K46").Value = 654.2308
$ws.Range("L46").Value = 1650
$ws.Range("M46").Value = -466.2308
$ws.Range("N46").Value = -2026

# LTW row 110
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 35514.4
$ws.Range("J110").Value = 35514.4
$ws.Range("L110").Value = 35514.4
$ws.Range("N110").Value = -43694.4

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3511.1538
$ws.Range("I126").Value = 3714.1052
$ws.Range("J126").Value = 2960.2856
$ws.Range("K126").Value = 11142.3156
$ws.Range("L126").Value = 8880.856800000001
$ws.Range("M126").Value = -8672.3156
$ws.Range("N126").Value = -13820.8568

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1985.7142
$ws.Range("I81").Value = 1816.6666
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3633.3332
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2572.3332
$ws.Range("N81").Value = -8122

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1985.7142
$ws.Range("I84").Value = 1816.6666
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 18166.666
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -12862.666
$ws.Range("N84").Value = -40608

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1059.6666
$ws.Range("I100").Value = 677.4286
$ws.Range("K100").Value = 1354.8572
$ws.Range("M100").Value = -813.8571999999999

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1522.1111
$ws.Range("I122").Value = 1600
$ws.Range("J122").Value = 1249.5
$ws.Range("K122").Value = 4800
$ws.Range("L122").Value = 3748.5
$ws.Range("M122").Value = -2350
$ws.Range("N122").Value = -8648.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1295.4524
$ws.Range("I132").Value = 1099.4814
$ws.Range("J132").Value = 1648.2
$ws.Range("K132").Value = 3298.4442
$ws.Range("L132").Value = 4944.6
$ws.Range("M132").Value = -768.4441999999999
$ws.Range("N132").Value = -10004.6

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2189.3333
$ws.Range("I136").Value = 800
$ws.Range("K136").Value = 2400
$ws.Range("M136").Value = 150
